# Supplemental_Tables_PB.xlsx edit:
#  - Swap the Mean/Median values that were originally entered in the
#    wrong columns on the "TableS5" sheet (rows 16 and 22).
#  - Update the selected cell / active sheet bookkeeping to match the
#    state the workbook was left in after the edit (TableS5 becomes the
#    active / tabSelected sheet, "Table S3" loses tabSelected).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Fix the swapped Mean / Median values on TableS5
# ---------------------------------------------------------------------
$wsS5 = $wb.Worksheets.Item("TableS5")

$b16 = $wsS5.Range("B16").Value2
$c16 = $wsS5.Range("C16").Value2
$wsS5.Range("B16").Value2 = $c16
$wsS5.Range("C16").Value2 = $b16

$b22 = $wsS5.Range("B22").Value2
$c22 = $wsS5.Range("C22").Value2
$wsS5.Range("B22").Value2 = $c22
$wsS5.Range("C22").Value2 = $b22

# ---------------------------------------------------------------------
# 2. Update selections on each sheet
# ---------------------------------------------------------------------

# "Table S3" is no longer the tab that is selected; the last selected
# cell there becomes G12.
$wsS3 = $wb.Worksheets.Item("Table S3")
[void]$wsS3.Activate()
[void]$wsS3.Range("G12").Select()

# TableS5 ends up being the active / tabSelected sheet, with F7 selected.
[void]$wsS5.Activate()
[void]$wsS5.Range("F7").Select()

# ---------------------------------------------------------------------
# 3. Best-effort restore of the workbook window geometry
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 18280
$win.Top = 0
$win.Width = 19060
$win.Height = 20560
$win.TabRatio = 500
